$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-6, columns A-J), row 6 is a newly added row.
$data = @(
    @(1, 4, 7, 7, 4, 3, -3, 43, 5),
    @(2, 0, 5, 5, 4, 5, -1, 65, 5),
    @(3, 3, 6, 4, 1, 1, -5, 21, 5),
    @(4, 0, 7, 4, 5, 4, -2, 54, 5),
    @(5, 4, 6, 6, 2, 2, -4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    # Column J keeps the shared string "train_dim2_1" on every data row.
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select()
